# Update cryptos list cell values per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.665.95"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "2.731.57"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'562.36"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").Value = "'159.24"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.23%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("E10").Value = "  +3.91%  "
$ws.Range("D11").Value = "'5.61"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.27%  "
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "3.214.82"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "'26.89"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.42%  "
$ws.Range("D15").Value = "63.523.80"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "2.736.00"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("E18").Value = "  +3.05%  "
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("D20").Value = "'353.28"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("E23").Value = "  -2.15%  "
$ws.Range("D24").Value = "'64.14"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.69%  "
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").Value = "'8.36"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("D28").Value = "0.0₃0902"
$ws.Range("E28").Value = "  +2.91%  "
$ws.Range("E29").Value = "  +2.32%  "
$ws.Range("D30").Value = "'7.19"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.82%  "
$ws.Range("E31").Value = "  +11.99%  "
$ws.Range("D32").Value = "'164.00"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.64%  "
$ws.Range("D33").Value = "'20.03"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  +3.63%  "
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("D38").Value = "'0.973"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("D39").Value = "'343.70"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +7.04%  "
$ws.Range("D40").Value = "'6.23"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.75%  "
$ws.Range("D41").Value = "'4.09"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").Value = "'38.37"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("D43").Value = "'21.78"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.38%  "
$ws.Range("D44").Value = "'21.01"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.621"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'133.69"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.0999"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0249"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").Value = "'11.07"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.12%  "
